# Rename the worksheet to reflect the unified "DataNode" concept
# (Property1 -> DataNode, per commit: "unify the conception of DataNode, DataTable, Entity").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Move/save the active selection to D42, matching where the author left the
# cursor when the sheet was last saved.
$ws.Range("D42").Select()
